# Update the cryptos worksheet (Sheet1) with refreshed price / volume data
# and a handful of re-ranked rows (14/15, 30-33, 48/49), matching the
# "Updated cryptos list" GitHub Actions commit.
#
# Note: several Price values (column D) are plain decimal-looking strings
# (e.g. "243.71"). Those are forced to stay as text (not auto-converted to
# numbers) with a leading apostrophe, and the style is then reset back to
# "Normal" so no stray number-format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.527.40"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.883.51"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'243.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.83%  "
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'42.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.27%  "
$ws.Range("D9").Value = "'0.334"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D11").Value = "'0.0994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "2.153.79"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "'12.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.68%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.902.12"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.687"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "35.493.49"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "'71.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").Value = "'243.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'12.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'170.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "'2.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.70%  "
$ws.Range("D27").Value = "'8.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.87%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").Value = "'0.126"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'0.952"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +26.76%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0565"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("B33").Value = "BinanceUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D33").Value = "'1.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "'4.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("D35").Value = "'1.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.87%  "
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("E37").Value = "  +9.41%  "
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("D39").Value = "'0.0205"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").Value = "'90.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "1.354.88"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").Value = "'15.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.07%  "
$ws.Range("D43").Value = "'13.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +48.31%  "
$ws.Range("D44").Value = "'0.0592"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.54%  "
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  +5.86%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'45.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +33.27%  "
$ws.Range("D50").Value = "2.071.11"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").Value = "'0.0691"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
